$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set date format for column D (row 32 is new; ensure consistent date style)
$ws.Range("D2:D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Ensure static columns are correct for the newly appended row 32
$ws.Range("A32").Value = 4
$ws.Range("B32").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C32").Value = "Los Lagos"
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = 300000000
$ws.Range("G32").Value = "Espárragos"
$ws.Range("I32").Value = "Primera"
$ws.Range("N32").Value = '$/kilo'
$ws.Range("O32").Value = "Provincia de Linares"
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = "Hortaliza"

# Updated weekly price/volume data (rows 2-32), columns D,H,J,K,L,M,P
$Dvals = @(44165, 44488, 44511, 44503, 44490, 44476, 44484, 44179, 44473, 44159, 44497, 44482, 44516, 44168, 44475, 44162, 44496, 44498, 44487, 44509, 44494, 44481, 44169, 44474, 44176, 44504, 44491, 44166, 44495, 44161, 44517)
$Hvals = @("Verde", "Sin especificar", "Sin especificar", "Sin especificar", "Sin especificar", "Sin especificar", "Sin especificar", "Verde", "Sin especificar", "Verde", "Sin especificar", "Sin especificar", "Sin especificar", "Verde", "Sin especificar", "Verde", "Sin especificar", "Sin especificar", "Sin especificar", "Sin especificar", "Sin especificar", "Sin especificar", "Verde", "Sin especificar", "Verde", "Sin especificar", "Sin especificar", "Verde", "Sin especificar", "Verde", "Sin especificar")
$Jvals = @(300, 600, 144, 72, 72, 60, 550, 200, 100, 600, 120, 72, 360, 200, 100, 700, 84, 500, 120, 550, 120, 300, 600, 250, 700, 180, 500, 500, 520, 300, 120)
$Kvals = @(1600, 1700, 1700, 1600, 1700, 2000, 1700, 1600, 2000, 1600, 1800, 2000, 1600, 1600, 2000, 1600, 1800, 1600, 1800, 1700, 1700, 1700, 1600, 2000, 1600, 1600, 1700, 1600, 1800, 1700, 1600)
$Lvals = @(1600, 1800, 1700, 1600, 1700, 2000, 1700, 1600, 2000, 1700, 1800, 2000, 1600, 1600, 2000, 1600, 1800, 1600, 1800, 1700, 1700, 2000, 1600, 2000, 1600, 1600, 1700, 1600, 1800, 1700, 1600)
$Mvals = @(1600, 1750, 1700, 1600, 1700, 2000, 1700, 1600, 2000, 1650, 1800, 2000, 1600, 1600, 2000, 1600, 1800, 1600, 1800, 1700, 1700, 1850, 1600, 2000, 1600, 1600, 1700, 1600, 1800, 1700, 1600)
$Pvals = @(1600, 1750, 1700, 1600, 1700, 2000, 1700, 1600, 2000, 1650, 1800, 2000, 1600, 1600, 2000, 1600, 1800, 1600, 1800, 1700, 1700, 1850, 1600, 2000, 1600, 1600, 1700, 1600, 1800, 1700, 1600)

for ($i = 0; $i -lt 31; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 4).Value = $Dvals[$i]
    $ws.Cells.Item($r, 8).Value = $Hvals[$i]
    $ws.Cells.Item($r, 10).Value = $Jvals[$i]
    $ws.Cells.Item($r, 11).Value = $Kvals[$i]
    $ws.Cells.Item($r, 12).Value = $Lvals[$i]
    $ws.Cells.Item($r, 13).Value = $Mvals[$i]
    $ws.Cells.Item($r, 16).Value = $Pvals[$i]
}
